$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) - update "想去人数" (want-to-go count) column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 298
$wsExhibit.Range("F3").Value = 1197
$wsExhibit.Range("F4").Value = 2667
$wsExhibit.Range("F5").Value = 236

# Sheet "全部类型" (all types) - same events, update matching rows in column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 298
$wsAll.Range("F5").Value = 1197
$wsAll.Range("F6").Value = 2667
$wsAll.Range("F8").Value = 236
